# Auto-generated edit script: updates crypto price/volume table (rows 2-51)
# matching the GitHub Actions data refresh described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.814.91'
$ws.Cells.Item(2, 5).Value = '  +3.89%  '
$ws.Cells.Item(3, 4).Value = '2.648.78'
$ws.Cells.Item(3, 5).Value = '  +1.32%  '
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
$ws.Cells.Item(5, 4).Value = '''567.87'
$ws.Cells.Item(5, 5).Value = '  +6.27%  '
$ws.Cells.Item(6, 4).Value = '''146.59'
$ws.Cells.Item(6, 5).Value = '  +2.81%  '
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 5).Value = '  +5.92%  '
$ws.Cells.Item(9, 4).Value = '2.660.84'
$ws.Cells.Item(9, 5).Value = '  +1.63%  '
$ws.Cells.Item(10, 5).Value = '  +0.21%  '
$ws.Cells.Item(11, 5).Value = '  +4.69%  '
$ws.Cells.Item(12, 5).Value = '  +6.84%  '
$ws.Cells.Item(13, 4).Value = '''0.344'
$ws.Cells.Item(13, 5).Value = '  +3.52%  '
$ws.Cells.Item(14, 4).Value = '3.117.29'
$ws.Cells.Item(14, 5).Value = '  +1.29%  '
$ws.Cells.Item(15, 4).Value = '60.707.32'
$ws.Cells.Item(15, 5).Value = '  +3.75%  '
$ws.Cells.Item(16, 4).Value = '''21.93'
$ws.Cells.Item(16, 5).Value = '  +5.56%  '
$ws.Cells.Item(17, 4).Value = '''0.0000138'
$ws.Cells.Item(17, 5).Value = '  +4.45%  '
$ws.Cells.Item(18, 4).Value = '2.667.09'
$ws.Cells.Item(18, 5).Value = '  +1.84%  '
$ws.Cells.Item(19, 4).Value = '''4.55'
$ws.Cells.Item(19, 5).Value = '  +3.14%  '
$ws.Cells.Item(20, 4).Value = '''342.71'
$ws.Cells.Item(20, 5).Value = '  +2.47%  '
$ws.Cells.Item(21, 4).Value = '''10.44'
$ws.Cells.Item(21, 5).Value = '  +3.09%  '
$ws.Cells.Item(22, 5).Value = '  +2.61%  '
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).Value = '''1.00'
$ws.Cells.Item(23, 5).Value = '  +0.11%  '
$ws.Cells.Item(24, 2).Value = 'Litecoin'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(24, 4).Value = '''66.77'
$ws.Cells.Item(24, 5).Value = '  +0.12%  '
$ws.Cells.Item(25, 2).Value = 'Polygon'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(25, 4).Value = '''0.442'
$ws.Cells.Item(25, 5).Value = '  +5.32%  '
$ws.Cells.Item(26, 2).Value = 'Kaspa'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(26, 4).Value = '''0.164'
$ws.Cells.Item(26, 5).Value = '  +1.88%  '
$ws.Cells.Item(27, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(27, 4).Value = '''0.998'
$ws.Cells.Item(27, 5).Value = '  -0.08%  '
$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).Value = '''7.39'
$ws.Cells.Item(28, 5).Value = '  +4.31%  '
$ws.Cells.Item(29, 2).Value = 'PEPE'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(29, 4).Value = '0.0₃0802'
$ws.Cells.Item(29, 5).Value = '  +9.47%  '
$ws.Cells.Item(30, 2).Value = 'USDe'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(30, 4).Value = '''0.999'
$ws.Cells.Item(30, 5).Value = '  +0.04%  '
$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(31, 4).Value = '''1.71'
$ws.Cells.Item(31, 5).Value = '  +4.60%  '
$ws.Cells.Item(32, 2).Value = 'Aptos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(32, 4).Value = '''6.22'
$ws.Cells.Item(32, 5).Value = '  +3.81%  '
$ws.Cells.Item(33, 2).Value = 'Monero'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(33, 4).Value = '''158.43'
$ws.Cells.Item(33, 5).Value = '  +3.29%  '
$ws.Cells.Item(34, 2).Value = 'EthereumClassic'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(34, 4).Value = '''19.24'
$ws.Cells.Item(34, 5).Value = '  +1.58%  '
$ws.Cells.Item(35, 2).Value = 'NEARProtocol'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(35, 4).Value = '''4.09'
$ws.Cells.Item(35, 5).Value = '  +4.78%  '
$ws.Cells.Item(36, 2).Value = 'SuiNetwork'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(36, 4).Value = '''0.911'
$ws.Cells.Item(36, 5).Value = '  +9.58%  '
$ws.Cells.Item(37, 2).Value = 'Fetch.AI'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(37, 4).Value = '''0.903'
$ws.Cells.Item(37, 5).Value = '  +10.47%  '
$ws.Cells.Item(38, 5).Value = '  +5.42%  '
$ws.Cells.Item(39, 2).Value = 'OKB'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(39, 4).Value = '''37.40'
$ws.Cells.Item(39, 5).Value = '  +0.95%  '
$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).Value = '''1.51'
$ws.Cells.Item(40, 5).Value = '  +6.72%  '
$ws.Cells.Item(41, 2).Value = 'Bittensor'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(41, 4).Value = '''303.76'
$ws.Cells.Item(41, 5).Value = '  +6.67%  '
$ws.Cells.Item(42, 2).Value = 'Filecoin'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(42, 4).Value = '''3.65'
$ws.Cells.Item(42, 5).Value = '  +1.86%  '
$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 4).Value = '''0.997'
$ws.Cells.Item(43, 5).Value = '  -0.25%  '
$ws.Cells.Item(44, 2).Value = 'Stellar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(44, 4).Value = '''0.0988'
$ws.Cells.Item(44, 5).Value = '  +4.78%  '
$ws.Cells.Item(45, 2).Value = 'Mantle'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(45, 4).Value = '''0.605'
$ws.Cells.Item(45, 5).Value = '  +1.63%  '
$ws.Cells.Item(46, 2).Value = 'Aave'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(46, 4).Value = '''129.25'
$ws.Cells.Item(46, 5).Value = '  +14.36%  '
$ws.Cells.Item(47, 4).Value = '''0.0545'
$ws.Cells.Item(47, 5).Value = '  +3.13%  '
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).Value = '''19.36'
$ws.Cells.Item(48, 5).Value = '  +1.80%  '
$ws.Cells.Item(49, 4).Value = '''10.70'
$ws.Cells.Item(49, 5).Value = '  +0.04%  '
$ws.Cells.Item(50, 2).Value = 'VeChain'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(50, 4).Value = '''0.0238'
$ws.Cells.Item(50, 5).Value = '  +5.40%  '
$ws.Cells.Item(51, 2).Value = 'RenderToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(51, 4).Value = '''4.65'
$ws.Cells.Item(51, 5).Value = '  +4.32%  '
